# chore: update Sheets via scheduled runner
#
# Refreshes cached market-board figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the per-class leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Pure data values - no formulas
# are used on these sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 251.25
$ws.Range("I8").Value = 251.25
$ws.Range("K8").Value = 753.75
$ws.Range("M8").Value = -614.75
$ws.Range("H40").Value = 960
$ws.Range("I40").Value = 925.7143
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 925.7143
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -750.7143
$ws.Range("N40").Value = -1550
$ws.Range("H51").Value = 3805.5557
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3968.75
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3968.75
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -4936.75
$ws.Range("H100").Value = 1668.8948
$ws.Range("I100").Value = 1550.6428
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1550.6428
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1009.6428
$ws.Range("N100").Value = -3082
$ws.Range("H118").Value = 1020.3158
$ws.Range("I118").Value = 632.3077
$ws.Range("J118").Value = 1861
$ws.Range("K118").Value = 1896.9231
$ws.Range("L118").Value = 5583
$ws.Range("M118").Value = -239.9231
$ws.Range("N118").Value = -8897
$ws.Range("H137").Value = 5055153.5
$ws.Range("I137").Value = 7942470
$ws.Range("J137").Value = 2349.9167
$ws.Range("K137").Value = 23827410
$ws.Range("L137").Value = 7049.750100000001
$ws.Range("M137").Value = -23824860
$ws.Range("N137").Value = -12149.7501
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2675903.8
$ws.Range("I2").Value = 2464
$ws.Range("J2").Value = 14706382
$ws.Range("K2").Value = 2464
$ws.Range("L2").Value = 14706382
$ws.Range("M2").Value = -2351
$ws.Range("N2").Value = -14706608
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H45").Value = 51487.6
$ws.Range("I45").Value = 67942.53
$ws.Range("K45").Value = 67942.53
$ws.Range("M45").Value = -67565.53
$ws.Range("H116").Value = 2675903.8
$ws.Range("I116").Value = 2464
$ws.Range("J116").Value = 14706382
$ws.Range("K116").Value = 2464
$ws.Range("L116").Value = 14706382
$ws.Range("M116").Value = -170
$ws.Range("N116").Value = -14710970
$ws.Range("H122").Value = 2289.4783
$ws.Range("I122").Value = 2029.9166
$ws.Range("J122").Value = 2572.6365
$ws.Range("K122").Value = 6089.7498
$ws.Range("L122").Value = 7717.9095
$ws.Range("M122").Value = -3639.7498
$ws.Range("N122").Value = -12617.9095
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2675903.8
$ws.Range("I3").Value = 2464
$ws.Range("J3").Value = 14706382
$ws.Range("K3").Value = 2464
$ws.Range("L3").Value = 14706382
$ws.Range("M3").Value = -2350
$ws.Range("N3").Value = -14706610
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 3000
$ws.Range("I15").Value = 3000
$ws.Range("K15").Value = 3000
$ws.Range("M15").Value = -2830
$ws.Range("H22").Value = 853
$ws.Range("I22").Value = 893.8125
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 893.8125
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -543.8125
$ws.Range("N22").Value = -900
$ws.Range("H31").Value = 2719.8333
$ws.Range("I31").Value = 2260
$ws.Range("J31").Value = 2949.75
$ws.Range("K31").Value = 2260
$ws.Range("L31").Value = 2949.75
$ws.Range("M31").Value = -1965
$ws.Range("N31").Value = -3539.75
$ws.Range("H34").Value = 2719.8333
$ws.Range("I34").Value = 2260
$ws.Range("J34").Value = 2949.75
$ws.Range("K34").Value = 2260
$ws.Range("L34").Value = 2949.75
$ws.Range("M34").Value = -2058
$ws.Range("N34").Value = -3353.75
$ws.Range("H94").Value = 929.2381
$ws.Range("J94").Value = 883.44446
$ws.Range("L94").Value = 883.44446
$ws.Range("N94").Value = -1785.44446
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("H132").Value = 3704.6072
$ws.Range("I132").Value = 3569.3684
$ws.Range("J132").Value = 3990.111
$ws.Range("K132").Value = 10708.1052
$ws.Range("L132").Value = 11970.333
$ws.Range("M132").Value = -8178.1052
$ws.Range("N132").Value = -17030.333
$ws.Range("H134").Value = 2005.1177
$ws.Range("I134").Value = 2313.1428
$ws.Range("J134").Value = 1507.5385
$ws.Range("K134").Value = 6939.428400000001
$ws.Range("L134").Value = 4522.6155
$ws.Range("M134").Value = -4404.428400000001
$ws.Range("N134").Value = -9592.6155
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 242.08333
$ws.Range("I92").Value = 300.25
$ws.Range("J92").Value = 125.75
$ws.Range("K92").Value = 900.75
$ws.Range("L92").Value = 377.25
$ws.Range("M92").Value = 347.25
$ws.Range("N92").Value = -2873.25
$ws.Range("H103").Value = 17000500
$ws.Range("I103").Value = 17000500
$ws.Range("K103").Value = 51001500
$ws.Range("M103").Value = -51000621
$ws.Range("H114").Value = 1427.92
$ws.Range("I114").Value = 729.2
$ws.Range("J114").Value = 2476
$ws.Range("K114").Value = 2187.6
$ws.Range("L114").Value = 7428
$ws.Range("M114").Value = 1066.4
$ws.Range("N114").Value = -13936
$ws.Range("H131").Value = 1390702.1
$ws.Range("I131").Value = 1781.3334
$ws.Range("J131").Value = 1516967.6
$ws.Range("K131").Value = 5344.0002
$ws.Range("L131").Value = 4550902.800000001
$ws.Range("M131").Value = -304.0002000000004
$ws.Range("N131").Value = -4560982.800000001
$ws.Range("H139").Value = 17244400
$ws.Range("I139").Value = 38463044
$ws.Range("J139").Value = 4250
$ws.Range("K139").Value = 115389132
$ws.Range("L139").Value = 12750
$ws.Range("M139").Value = -115383992
$ws.Range("N139").Value = -23030
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 4751.25
$ws.Range("I6").Value = 4751.25
$ws.Range("K6").Value = 4751.25
$ws.Range("M6").Value = -4638.25
$ws.Range("H16").Value = 4751.25
$ws.Range("I16").Value = 4751.25
$ws.Range("K16").Value = 4751.25
$ws.Range("M16").Value = -4501.25
$ws.Range("H122").Value = 4720.077
$ws.Range("I122").Value = 6049.467
$ws.Range("J122").Value = 2907.2727
$ws.Range("K122").Value = 18148.401
$ws.Range("L122").Value = 8721.8181
$ws.Range("M122").Value = -15698.401
$ws.Range("N122").Value = -13621.8181
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 34731.332
$ws.Range("I122").Value = 34731.332
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 104193.996
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -101743.996
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4045.4822
$ws.Range("I132").Value = 4353.6445
$ws.Range("J132").Value = 2784.818
$ws.Range("K132").Value = 13060.9335
$ws.Range("L132").Value = 8354.454000000002
$ws.Range("M132").Value = -10530.9335
$ws.Range("N132").Value = -13414.454
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 18201.143
$ws.Range("J12").Value = 18201.143
$ws.Range("L12").Value = 18201.143
$ws.Range("N12").Value = -18485.143
$ws.Range("H132").Value = 5229.5776
$ws.Range("I132").Value = 6029.794
$ws.Range("J132").Value = 2756.182
$ws.Range("K132").Value = 18089.382
$ws.Range("L132").Value = 8268.545999999998
$ws.Range("M132").Value = -15559.382
$ws.Range("N132").Value = -13328.546
$ws.Range("H136").Value = 30243.732
$ws.Range("I136").Value = 7286.3
$ws.Range("J136").Value = 92854.91
$ws.Range("K136").Value = 21858.9
$ws.Range("L136").Value = 278564.73
$ws.Range("M136").Value = -19308.9
$ws.Range("N136").Value = -283664.73
